$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# A leading apostrophe is used for numeric-looking Price values so Excel
# keeps them stored as text (matching the source data), rather than
# auto-converting them to numbers.

$ws.Range("D2").Value = "26.367.40"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "1.794.40"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'306.93"
$ws.Range("D7").Value = "'0.4529"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("D8").Value = "'0.3593"
$ws.Range("E8").Value = "  -2.23%  "
$ws.Range("D9").Value = "'46.06"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "'0.07076"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "'0.8855"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D13").Value = "'19.44"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "1.808.86"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'5.277"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "'85.02"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "'1.005"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'0.000008500"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "26.394.04"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").Value = "'4.959"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "2.042.87"
$ws.Range("D25").Value = "'10.53"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'1.967"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("D27").Value = "'150.97"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'17.79"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").Value = "'2.024"
$ws.Range("E29").Value = "  +3.16%  "
$ws.Range("D30").Value = "'111.82"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").Value = "'4.849"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "'0.08690"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("D33").Value = "'3.117"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").Value = "'2.793"
$ws.Range("E34").Value = "  +11.61%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").Value = "'0.7203"
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("D37").Value = "'1.103"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("D38").Value = "'1.002"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "'0.05066"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "'2.853"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").Value = "'0.5082"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Value = "'6.816"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").Value = "'0.1514"
$ws.Range("E45").Value = "  -5.11%  "
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "'0.4629"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").Value = "'9.870"
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").Value = "'100.93"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  -2.35%  "
